{"js": "// The document ends with a \"Requisitos\" section that lists\n// \"8800009: Canto Coral I (Requisito)\" followed by a blank paragraph and\n// two footer-style paragraphs (\"Ver no Jupiter Salvar em pdf Salvar em\n// docx\" and the \"\u00a9 2020 ...\" copyright line). This edit removes that\n// trailing blank paragraph and the two footer paragraphs, so the\n// requirement line is immediately followed by the final blank paragraph\n// that precedes the page-break paragraph.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the paragraphs to remove by their exact text content so the\n// edit is robust to the surrounding document rather than relying on a\n// fixed index.\nconst requisitoText = \"8800009: Canto Coral I (Requisito)\";\nconst jupiterText = \"Ver no Jupiter Salvar em pdf Salvar em docx\";\nconst copyrightText =\n  \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\";\n\nlet requisitoIdx = -1;\nlet jupiterIdx = -1;\nlet copyrightIdx = -1;\n\nfor (let i = 0; i < items.length; i++) {\n  const text = items[i].text;\n  if (requisitoIdx === -1 && text === requisitoText) {\n    requisitoIdx = i;\n  } else if (jupiterIdx === -1 && text === jupiterText) {\n    jupiterIdx = i;\n  } else if (copyrightIdx === -1 && text === copyrightText) {\n    copyrightIdx = i;\n  }\n}\n\n// The blank paragraph sits right after the requirement paragraph, and\n// immediately before the \"Ver no Jupiter...\" paragraph.\nconst blankIdx = requisitoIdx >= 0 ? requisitoIdx + 1 : -1;\n\n// Delete bottom-up so earlier indices stay valid.\nif (copyrightIdx !== -1) {\n  items[copyrightIdx].delete();\n}\nif (jupiterIdx !== -1) {\n  items[jupiterIdx].delete();\n}\nif (blankIdx !== -1 && blankIdx === jupiterIdx - 1) {\n  items[blankIdx].delete();\n}\n\nawait context.sync();\n", "ps1": "# The document ends with a \"Requisitos\" section that lists\n# \"8800009: Canto Coral I (Requisito)\" followed by a blank paragraph and\n# two footer-style paragraphs (\"Ver no Jupiter Salvar em pdf Salvar em\n# docx\" and the \"(c) 2020 ...\" copyright line). This edit removes that\n# trailing blank paragraph and the two footer paragraphs, so the\n# requirement line is immediately followed by the final blank paragraph\n# that precedes the page-break paragraph.\n\n$d = $word.ActiveDocument\n\n$requisitoText = \"8800009: Canto Coral I (Requisito)\"\n$jupiterText = \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n$copyrightText = [char]0x00A9 + \" 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n\n# Locate the paragraphs to remove by their exact text content so the\n# edit is robust to the surrounding document rather than relying on a\n# fixed index.\n$count = $d.Paragraphs.Count\n$requisitoIndex = -1\n$jupiterIndex = -1\n$copyrightIndex = -1\n\nfor ($i = 1; $i -le $count; $i++) {\n    $txt = $d.Paragraphs($i).Range.Text.TrimEnd(\"`r\")\n    if ($requisitoIndex -eq -1 -and $txt -eq $requisitoText) {\n        $requisitoIndex = $i\n    } elseif ($jupiterIndex -eq -1 -and $txt -eq $jupiterText) {\n        $jupiterIndex = $i\n    } elseif ($copyrightIndex -eq -1 -and $txt -eq $copyrightText) {\n        $copyrightIndex = $i\n    }\n}\n\n# The blank paragraph sits right after the requirement paragraph, and\n# immediately before the \"Ver no Jupiter...\" paragraph.\n$blankIndex = -1\nif ($requisitoIndex -ne -1) {\n    $blankIndex = $requisitoIndex + 1\n}\n\n# Delete bottom-up (highest index first) so earlier indices stay valid.\nif ($copyrightIndex -ne -1) {\n    $d.Paragraphs($copyrightIndex).Range.Delete()\n}\nif ($jupiterIndex -ne -1) {\n    $d.Paragraphs($jupiterIndex).Range.Delete()\n}\nif ($blankIndex -ne -1 -and $blankIndex -eq ($jupiterIndex - 1)) {\n    $d.Paragraphs($blankIndex).Range.Delete()\n}\n"}
